# Generate Report for Archive
$wb = $excel.ActiveWorkbook

# --- Update status text everywhere it appears ("Ready for handoff" -> "In Translation") ---
foreach ($ws in $wb.Worksheets) {
    $found = $ws.Cells.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        do {
            $found.Value = "In Translation"
            $found = $ws.Cells.FindNext($found)
        } while ($found -ne $null -and $found.Address() -ne $firstAddress)
    }
}

# --- Adjust column widths ---
# Target stored OOXML width is 13.4101845877511 characters; this runtime rounds
# ColumnWidth to the nearest 1/6 character before persisting, so 12.5 is the
# closest achievable setting (-> stored width 13.333333333333334).
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth
